# Generate Report for Handoff
# Adds two new source files (two .png files) to the localization status
# report, alongside the existing .md file, across the Overview, zh-cn and
# de-de worksheets. Also refreshes the "latest handoff" timestamps and the
# GUID used for the primary markdown source file.

$wb = $excel.ActiveWorkbook

$oldGuid = "2c179622-9275-43d8-a71d-7962e3f1b746"
$newGuid = "1c51fda2-ba03-41cd-aed2-1d5f364f33cd"

$mdFile  = "$newGuid.md"
$png1    = "41a28e61-2390-490c-9300-1306ffb75006.png"
$png2    = "f9225f2d-a267-4b62-9386-4e9a55bdcf30.png"

$zhHash  = "b0826e70d58049e683eafe8c6ffdc1d94e3f0d8d"
$deHash  = "b0826e70d58049e683eafe8c6ffdc1d94e3f0d8d"

$zhTarget = "$newGuid.$zhHash.zh-cn.xlf"
$deTarget = "$newGuid.$deHash.de-de.xlf"

$png1Target = "ff8efbe0af6421feee92992509fabbb53a46621d.png"
$png2Target = "6b479b4f212dc613b1ec136ed7fabcd14af1b71a.png"

$status       = "Ready for handoff"
$overviewDate = "2016-01-13 13:01:02"
$zhDate       = "2016-03-13 13:00:58"
$deDate       = "2016-03-13 13:01:02"
$epoch        = "0001-01-01 00:00:00"

$srcCommit    = "d72d400e2391c3bd2d479448e17e7e0f4e95317c"
$zhCommit     = "1527d2413ccc8c21e09362fe775baa1ebbbabf52"
$deCommit     = "7c09eb37b911c417a84a1fb3e76a491340aad54d"

function SrcUrl($name) {
    return "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$name"
}
function ZhHandoffUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$name"
}
function DeHandoffUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$name"
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Hyperlinks.Delete()

$wsOv.Range("A2").Value = $mdFile
$wsOv.Range("B2").Value = $status
$wsOv.Range("C2").Value = $status
$wsOv.Range("D2").Value = $overviewDate

$wsOv.Range("A3").Value = $png1
$wsOv.Range("B3").Value = $status
$wsOv.Range("C3").Value = $status
$wsOv.Range("D3").Value = $overviewDate

$wsOv.Range("A4").Value = $png2
$wsOv.Range("B4").Value = $status
$wsOv.Range("C4").Value = $status
$wsOv.Range("D4").Value = $overviewDate

$wsOv.Hyperlinks.Add($wsOv.Range("A2"), (SrcUrl $mdFile), "", "", $mdFile) | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("A3"), (SrcUrl $png1), "", "", $png1) | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("A4"), (SrcUrl $png2), "", "", $png2) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()

$wsZh.Range("A2").Value = $mdFile
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = $status
$wsZh.Range("D2").Value = $zhTarget
$wsZh.Range("E2").Value = $zhDate
$wsZh.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H2").Value = $epoch
$wsZh.Range("I2").Value = "Include"

$wsZh.Range("A3").Value = $png1
$wsZh.Range("B3").Value = ".png"
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = $png1Target
$wsZh.Range("E3").Value = $zhDate
$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value = $epoch
$wsZh.Range("I3").Value = "IsDependency"
$wsZh.Range("J3").Value = "e2e\$mdFile"

$wsZh.Range("A4").Value = $png2
$wsZh.Range("B4").Value = ".png"
$wsZh.Range("C4").Value = $status
$wsZh.Range("D4").Value = $png2Target
$wsZh.Range("E4").Value = $zhDate
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H4").Value = $epoch
$wsZh.Range("I4").Value = "IsDependency"
$wsZh.Range("J4").Value = "e2e\$mdFile"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), (SrcUrl $mdFile), "", "", $mdFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), (SrcUrl $mdFile), "", "", ".md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), (ZhHandoffUrl $zhTarget), "", "", $zhTarget) | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), (SrcUrl $png1), "", "", $png1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), (SrcUrl $png1), "", "", ".png") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), (ZhHandoffUrl $png1Target), "", "", $png1Target) | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), (SrcUrl $png2), "", "", $png2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), (SrcUrl $png2), "", "", ".png") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), (ZhHandoffUrl $png2Target), "", "", $png2Target) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()

$wsDe.Range("A2").Value = $mdFile
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = $status
$wsDe.Range("D2").Value = $deTarget
$wsDe.Range("E2").Value = $deDate
$wsDe.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H2").Value = $epoch
$wsDe.Range("I2").Value = "Include"

$wsDe.Range("A3").Value = $png1
$wsDe.Range("B3").Value = ".png"
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = $png1Target
$wsDe.Range("E3").Value = $deDate
$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value = $epoch
$wsDe.Range("I3").Value = "IsDependency"
$wsDe.Range("J3").Value = "e2e\$mdFile"

$wsDe.Range("A4").Value = $png2
$wsDe.Range("B4").Value = ".png"
$wsDe.Range("C4").Value = $status
$wsDe.Range("D4").Value = $png2Target
$wsDe.Range("E4").Value = $deDate
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H4").Value = $epoch
$wsDe.Range("I4").Value = "IsDependency"
$wsDe.Range("J4").Value = "e2e\$mdFile"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), (SrcUrl $mdFile), "", "", $mdFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), (SrcUrl $mdFile), "", "", ".md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), (DeHandoffUrl $deTarget), "", "", $deTarget) | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), (SrcUrl $png1), "", "", $png1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), (SrcUrl $png1), "", "", ".png") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), (DeHandoffUrl $png1Target), "", "", $png1Target) | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), (SrcUrl $png2), "", "", $png2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), (SrcUrl $png2), "", "", ".png") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), (DeHandoffUrl $png2Target), "", "", $png2Target) | Out-Null

Write-Host "Report generated for handoff."
